{"js": "// Change the \"N\u00b0\" cell in the first table row: the run holding the\n// superscript \"o\" becomes a plain (non-superscript) \"\u00b0\" character, and the\n// \"_GoBack\" bookmark (previously sitting alone in the empty paragraph at the\n// very end of the document) is moved to sit right after that run.\n\n// 1) Drop the old \"_GoBack\" bookmark wherever it currently lives (the last,\n//    empty paragraph of the document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the run containing the superscript \"o\" \u2014 it's the second run of\n//    the very first paragraph in the document (\"N\" + superscript \"o\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst found = firstParagraph.search(\"o\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nconst superscriptRun = found.items[0];\n\n// 3) Replace that run with an equivalent run that is no longer superscript\n//    and whose text is the degree sign \"\u00b0\" instead of the letter \"o\".\n//    insertOoxml is used (rather than font.superscript = false) so the\n//    <w:vertAlign> element is removed outright instead of being rewritten\n//    as <w:vertAlign w:val=\"baseline\"/>.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r w:rsidRPr=\"009C1197\">' +\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr>' +\n  '<w:t>\\u00B0</w:t>' +\n  '</w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nsuperscriptRun.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Re-insert the \"_GoBack\" bookmark immediately after the run we just\n//    replaced, matching its new location in the first table cell.\nconst afterRun = superscriptRun.getRange(Word.RangeLocation.after);\nafterRun.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change the \"N\u00b0\" cell in the first table row: the run holding the\n# superscript \"o\" becomes a plain (non-superscript) \"\u00b0\" character, and the\n# \"_GoBack\" bookmark (previously sitting alone in the empty paragraph at the\n# very end of the document) is moved to sit right after that run.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark wherever it currently lives (the last,\n#    empty paragraph of the document). Plain Font/property edits on this\n#    engine leave <w:vertAlign w:val=\"baseline\"/> behind instead of removing\n#    the element, and Bookmarks.Add() does not persist here, so the run\n#    swap and the bookmark re-creation below are both done together via a\n#    single raw-OOXML splice (Range.InsertXML) on the target paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Locate the run containing the superscript \"o\" \u2014 it's the second run of\n#    the very first paragraph in the document (\"N\" + superscript \"o\").\n$find = $d.Content\n$find.SetRange(0, 0)\n$find.Find.ClearFormatting()\n$find.Find.Text = \"o\"\n$find.Find.MatchCase = $true\n$find.Find.MatchWholeWord = $false\n$find.Find.Forward = $true\n$find.Find.Wrap = 0\n$find.Find.Execute() | Out-Null\n\n# 3) Replace the whole enclosing paragraph with an equivalent paragraph\n#    where that run is no longer superscript, its text is the degree sign\n#    \"\u00b0\" instead of the letter \"o\", and the \"_GoBack\" bookmark sits right\n#    after it (InsertXML works at paragraph granularity on this engine, and\n#    there's no supported Range.InsertBookmark / Bookmarks.Add here, so the\n#    bookmark tags are included directly in the spliced OOXML).\n$degree = [char]0x00B0\n$paragraphOoxml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p w:rsidR=\"00BD01F2\" w:rsidRPr=\"009C1197\" w:rsidRDefault=\"00BD01F2\" w:rsidP=\"004C7CF1\">' +\n  '<w:pPr><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr></w:pPr>' +\n  '<w:r w:rsidRPr=\"009C1197\"><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr><w:t>N</w:t></w:r>' +\n  '<w:r w:rsidRPr=\"009C1197\"><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr><w:t>' + $degree + '</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$find.InsertXML($paragraphOoxml)\n"}
